$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for each coin row ---
$ws.Range("D2").Value = "30.439.32"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.877.11"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.20"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4807"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2824"
$ws.Range("E8").Value = "  -2.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06526"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").Value = "1.874.63"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07495"
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.58"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.075"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.57"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6629"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "30.376.88"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007605"
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("D20").Value = "2.111.26"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.307"
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.192"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.328"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.98"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.46"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.965"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.460"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09413"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.302"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.035"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05014"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.213"
$ws.Range("E34").Value = "  +5.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7436"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.704"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01828"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.616"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.43"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.835"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4273"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.445"
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.34"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1277"
$ws.Range("E47").Value = "  -6.53%  "
$ws.Range("E48").Value = "  -6.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.851"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.68"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3882"
$ws.Range("E51").Value = "  +0.53%  "

# --- Rows 22/23: BinanceUSD and BitcoinCash swapped positions ---
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "222.49"
$ws.Range("E22").Value = "  +16.53%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9987"
$ws.Range("E23").Value = "  -0.20%  "

# --- Rows 39/40: RenderToken and TrustWalletToken swapped positions ---
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9061"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.065"
$ws.Range("E40").Value = "  -0.38%  "
